# Generate Report for Handoff
# Adds two new handed-off files (8d1d55c9-...md and c15ac63d-...md) to all
# three worksheets: Overview, zh-cn, de-de. Each sheet's table grows by two
# rows (from 3 to 5 total rows including the header).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"  (table3.xml / sheet1.xml)  columns A-G
#   A File Name | B Path And Name | C Extension | D Publish URL
#   E zh-cn | F de-de | G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

# Row 4: 8d1d55c9-92f9-4f1a-bdd7-99adb74e440e.md
$wsOverview.Range("A4").Value = "8d1d55c9-92f9-4f1a-bdd7-99adb74e440e.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f9f90abc022d575668fb15501eb58fa66a598380/e2e/8d1d55c9-92f9-4f1a-bdd7-99adb74e440e.md", "", "", "e2e\8d1d55c9-92f9-4f1a-bdd7-99adb74e440e.md") | Out-Null
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-09-05 06:46:23"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 5: c15ac63d-29f9-4b84-8719-5cc5982dd3e1.md
$wsOverview.Range("A5").Value = "c15ac63d-29f9-4b84-8719-5cc5982dd3e1.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/69841d40e2694d9c86477eb7eb4a63344882ece6/e2e/c15ac63d-29f9-4b84-8719-5cc5982dd3e1.md", "", "", "e2e\c15ac63d-29f9-4b84-8719-5cc5982dd3e1.md") | Out-Null
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-09-05 06:46:23"
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "zh-cn"  (table1.xml / sheet2.xml)  columns A-P
#   A Source File Name | B File Extension | C Status | D Source Path
#   E Priority | F Content Duplicate | G Latest Handoff File
#   H Latest Handoff Datetime | I Latest Target File | J Latest Handback File
#   K Latest Handback DateTime | L Reference Tokens | M To be localized
#   N Dependency From | O Has metadata | P Error Detail
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null
$loZhCn.ListRows.Add() | Out-Null

# Row 4: 8d1d55c9-92f9-4f1a-bdd7-99adb74e440e.md
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f9f90abc022d575668fb15501eb58fa66a598380/e2e/8d1d55c9-92f9-4f1a-bdd7-99adb74e440e.md", "", "", "8d1d55c9-92f9-4f1a-bdd7-99adb74e440e.md") | Out-Null
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'False"
$wsZhCn.Range("G4").Value = "8d1d55c9-92f9-4f1a-bdd7-99adb74e440e.4eb0e07a36bf56f995024249af6f48bb3b14e18d.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-09-05 06:46:17"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("O4").Value = "'False"

# Row 5: c15ac63d-29f9-4b84-8719-5cc5982dd3e1.md
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/69841d40e2694d9c86477eb7eb4a63344882ece6/e2e/c15ac63d-29f9-4b84-8719-5cc5982dd3e1.md", "", "", "c15ac63d-29f9-4b84-8719-5cc5982dd3e1.md") | Out-Null
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("F5").Value = "'False"
$wsZhCn.Range("G5").Value = "c15ac63d-29f9-4b84-8719-5cc5982dd3e1.d5cf2fcc7ef680105d940f2d3ba6ad7e0620d671.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-09-05 06:46:17"
$wsZhCn.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M5").Value = "'True"
$wsZhCn.Range("O5").Value = "'False"

# ---------------------------------------------------------------------
# Sheet "de-de"  (table2.xml / sheet3.xml)  columns A-P  (same layout as zh-cn)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null
$loDeDe.ListRows.Add() | Out-Null

# Row 4: 8d1d55c9-92f9-4f1a-bdd7-99adb74e440e.md
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f9f90abc022d575668fb15501eb58fa66a598380/e2e/8d1d55c9-92f9-4f1a-bdd7-99adb74e440e.md", "", "", "8d1d55c9-92f9-4f1a-bdd7-99adb74e440e.md") | Out-Null
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'False"
$wsDeDe.Range("G4").Value = "8d1d55c9-92f9-4f1a-bdd7-99adb74e440e.4eb0e07a36bf56f995024249af6f48bb3b14e18d.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-09-05 06:46:23"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("O4").Value = "'False"

# Row 5: c15ac63d-29f9-4b84-8719-5cc5982dd3e1.md
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/69841d40e2694d9c86477eb7eb4a63344882ece6/e2e/c15ac63d-29f9-4b84-8719-5cc5982dd3e1.md", "", "", "c15ac63d-29f9-4b84-8719-5cc5982dd3e1.md") | Out-Null
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("F5").Value = "'False"
$wsDeDe.Range("G5").Value = "c15ac63d-29f9-4b84-8719-5cc5982dd3e1.d5cf2fcc7ef680105d940f2d3ba6ad7e0620d671.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-09-05 06:46:23"
$wsDeDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M5").Value = "'True"
$wsDeDe.Range("O5").Value = "'False"
